$d = $word.ActiveDocument
$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Edit 1 --------------------------------------------------------------
# "1차 시도 (강의자료 참고)" -> split the trailing run into three runs:
#   " (강의자료" / " 내림차순 정렬" / " 참고)"
$rng1 = $d.Content
$rng1.Find.Execute("1차 시도 (강의자료 참고)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $rng1.Paragraphs(1)
$paraRange1 = $para1.Range

$xml1 = "<w:p $wordNs>" + `
  "<w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr><w:t>1차 시도</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr><w:t xml:space='preserve'> (강의자료</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr><w:t xml:space='preserve'> 내림차순 정렬</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr><w:t xml:space='preserve'> 참고)</w:t></w:r>" + `
  "</w:p>"

$paraRange1.InsertXML($xml1)

# --- Edit 2 ----------------------------------------------------------------
# After "이렇게 하면 h값이 나오는 코드를 구현하였다." add three new paragraphs
# (one blank, a bold "실행결과 (실패)" heading, and a result description),
# keeping the existing trailing blank paragraph at the end of the document.
$rng2 = $d.Content
$rng2.Find.Execute("이렇게 하면 h값이 나오는 코드를 구현하였다.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1)
$lastPara = $para2.Next()
$lastParaRange = $lastPara.Range
$lastParaRange.Collapse(1)

$xml2 = "<w:p $wordNs><w:pPr><w:tabs><w:tab w:val='left' w:pos='1200'/><w:tab w:val='left' w:pos='6285'/></w:tabs></w:pPr></w:p>" + `
  "<w:p $wordNs><w:pPr><w:tabs><w:tab w:val='left' w:pos='1200'/><w:tab w:val='left' w:pos='6285'/></w:tabs><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:b/><w:bCs/></w:rPr><w:t>실행결과 (실패)</w:t></w:r></w:p>" + `
  "<w:p $wordNs><w:pPr><w:tabs><w:tab w:val='left' w:pos='1200'/><w:tab w:val='left' w:pos='6285'/></w:tabs></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t xml:space='preserve'>예상과 달리 citations = </w:t></w:r>" + `
  "<w:r><w:t>[3, 0, 6, 1, 5]</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t xml:space='preserve'> 상태에서 </w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>h값이 3이 아닌 0이 반환되었다.</w:t></w:r></w:p>" + `
  "<w:p $wordNs><w:pPr><w:tabs><w:tab w:val='left' w:pos='1200'/><w:tab w:val='left' w:pos='6285'/></w:tabs><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr></w:pPr></w:p>"

$lastParaRange.InsertXML($xml2)

Write-Output "edit complete"
